# Update latest output (run 104)
$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update Cost ($) and Unit Cost ($/ML) columns for rows 3-5 ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E3").Value = 52.77443925000003
$wsSchedule.Range("F3").Value = 1.469630722639934
$wsSchedule.Range("E4").Value = 457.1512432499999
$wsSchedule.Range("F4").Value = 30.23487058531746
$wsSchedule.Range("E5").Value = -59.26747124999999
$wsSchedule.Range("F5").Value = -1.742136133156966

# --- Sheet "Detailed": update Price column (and a couple of Type labels) ---
$wsDetailed = $wb.Worksheets.Item("Detailed")
$wsDetailed.Range("B29").Value = -5.00468
$wsDetailed.Range("B30").Value = -0.09135
$wsDetailed.Range("B31").Value = 0.00005
$wsDetailed.Range("B32").Value = 0.01929
$wsDetailed.Range("C32").Value = "historical"
$wsDetailed.Range("B33").Value = 0.51
$wsDetailed.Range("C33").Value = "historical"
$wsDetailed.Range("B34").Value = 0
$wsDetailed.Range("B35").Value = -0.92309
$wsDetailed.Range("B36").Value = 0.264
$wsDetailed.Range("B37").Value = 12.24467
$wsDetailed.Range("B38").Value = 26.68735
$wsDetailed.Range("B39").Value = 45.84936
$wsDetailed.Range("B40").Value = 59.11456
$wsDetailed.Range("B41").Value = 64.01260000000001
$wsDetailed.Range("B43").Value = 62.62134
$wsDetailed.Range("B44").Value = 64.05184
$wsDetailed.Range("B45").Value = 65
$wsDetailed.Range("B46").Value = 59.37719
$wsDetailed.Range("B47").Value = 57.96819
$wsDetailed.Range("B48").Value = 58.16113
$wsDetailed.Range("B49").Value = 63.63924
$wsDetailed.Range("B50").Value = 61.77255
$wsDetailed.Range("B55").Value = 57.06003
$wsDetailed.Range("B61").Value = 65.67543000000001
$wsDetailed.Range("B62").Value = 65.2897
$wsDetailed.Range("B63").Value = 61.31557
$wsDetailed.Range("B65").Value = 8.71917
$wsDetailed.Range("B66").Value = 0.7
$wsDetailed.Range("B67").Value = 0
$wsDetailed.Range("B68").Value = -2.535
$wsDetailed.Range("B69").Value = -5.68809
$wsDetailed.Range("B70").Value = -6.48267
$wsDetailed.Range("B71").Value = -7.87948
$wsDetailed.Range("B72").Value = -8.815630000000001
$wsDetailed.Range("B73").Value = -8.67817
$wsDetailed.Range("B74").Value = -9.99
$wsDetailed.Range("B75").Value = -9.99
$wsDetailed.Range("B76").Value = -8.3056
$wsDetailed.Range("B77").Value = -8.293850000000001
$wsDetailed.Range("B78").Value = -8.12302
$wsDetailed.Range("B79").Value = -7.98411
$wsDetailed.Range("B80").Value = -7.00778
$wsDetailed.Range("B81").Value = -6.49292
$wsDetailed.Range("B82").Value = -5.51
$wsDetailed.Range("B83").Value = -5.16677
$wsDetailed.Range("B84").Value = -2.9952
$wsDetailed.Range("B85").Value = 0.01019
$wsDetailed.Range("B86").Value = 9.523630000000001
$wsDetailed.Range("B87").Value = 55.62763
$wsDetailed.Range("B88").Value = 72.86584000000001
$wsDetailed.Range("B89").Value = 72.78848000000001
$wsDetailed.Range("B90").Value = 71.8766
$wsDetailed.Range("B91").Value = 65
$wsDetailed.Range("B92").Value = 59.15606
$wsDetailed.Range("B93").Value = 57.3
$wsDetailed.Range("B94").Value = 57.3
$wsDetailed.Range("B95").Value = 59.93919
$wsDetailed.Range("B96").Value = 63.27752
$wsDetailed.Range("B97").Value = 64.12223
